$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15 with value matching style of A14 (time format)
$ws.Range("A15").Value = 0.020231481481481482

# Update formulas to extend range to row 15
$ws.Range("C2").Formula = "=SUM(A2:A15)"
$ws.Range("B3").Formula = "=SUM(A9:A15)"

# Copy style from A14 to A15 so formatting matches
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to D6
$ws.Range("D6").Select()
